$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# Row 2
Set-TextValue $ws.Range("D2") '25.866.01'
Set-TextValue $ws.Range("E2") '  -0.16%  '

# Row 3
Set-TextValue $ws.Range("D3") '1.639.50'
Set-TextValue $ws.Range("E3") '  +0.26%  '

# Row 4
Set-TextValue $ws.Range("D4") '1.003'
Set-TextValue $ws.Range("E4") '  +0.14%  '

# Row 5
Set-TextValue $ws.Range("D5") '215.71'
Set-TextValue $ws.Range("E5") '  +0.84%  '

# Row 6
Set-TextValue $ws.Range("D6") '0.5066'
Set-TextValue $ws.Range("E6") '  +0.29%  '

# Row 7
Set-TextValue $ws.Range("E7") '  +0.19%  '

# Row 8
Set-TextValue $ws.Range("D8") '0.2582'
Set-TextValue $ws.Range("E8") '  +0.54%  '

# Row 9
Set-TextValue $ws.Range("D9") '0.06440'
Set-TextValue $ws.Range("E9") '  +1.39%  '

# Row 10
Set-TextValue $ws.Range("D10") '19.61'
Set-TextValue $ws.Range("E10") '  -0.23%  '

# Row 11
Set-TextValue $ws.Range("D11") '0.07784'

# Row 12
Set-TextValue $ws.Range("D12") '4.290'
Set-TextValue $ws.Range("E12") '  +0.44%  '

# Row 13
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue $ws.Range("D13") '1.641.23'
Set-TextValue $ws.Range("E13") '  +0.28%  '

# Row 14
$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextValue $ws.Range("D14") '1.865.67'
Set-TextValue $ws.Range("E14") '  +0.24%  '

# Row 15
Set-TextValue $ws.Range("D15") '0.5639'
Set-TextValue $ws.Range("E15") '  +3.90%  '

# Row 16
Set-TextValue $ws.Range("D16") '0.0₅7620'
Set-TextValue $ws.Range("E16") '  -1.34%  '

# Row 17
Set-TextValue $ws.Range("D17") '63.19'
Set-TextValue $ws.Range("E17") '  -1.30%  '

# Row 18
Set-TextValue $ws.Range("D18") '25.882.31'
Set-TextValue $ws.Range("E18") '  -0.19%  '

# Row 19
Set-TextValue $ws.Range("E19") '  +0.14%  '

# Row 20
Set-TextValue $ws.Range("D20") '195.30'
Set-TextValue $ws.Range("E20") '  +0.25%  '

# Row 21
Set-TextValue $ws.Range("D21") '4.325'
Set-TextValue $ws.Range("E21") '  -2.25%  '

# Row 22
Set-TextValue $ws.Range("D22") '9.890'
Set-TextValue $ws.Range("E22") '  -0.11%  '

# Row 23
Set-TextValue $ws.Range("D23") '6.104'
Set-TextValue $ws.Range("E23") '  +0.13%  '

# Row 25
Set-TextValue $ws.Range("D25") '1.796'
Set-TextValue $ws.Range("E25") '  -4.84%  '

# Row 26
Set-TextValue $ws.Range("D26") '0.1276'
Set-TextValue $ws.Range("E26") '  +3.17%  '

# Row 27
Set-TextValue $ws.Range("D27") '140.00'
Set-TextValue $ws.Range("E27") '  -1.96%  '

# Row 28
Set-TextValue $ws.Range("D28") '6.813'
Set-TextValue $ws.Range("E28") '  +0.15%  '

# Row 29
Set-TextValue $ws.Range("D29") '15.48'
Set-TextValue $ws.Range("E29") '  -0.47%  '

# Row 30
Set-TextValue $ws.Range("E30") '  +0.67%  '

# Row 31
Set-TextValue $ws.Range("D31") '0.04880'
Set-TextValue $ws.Range("E31") '  +0.32%  '

# Row 32
Set-TextValue $ws.Range("D32") '3.306'
Set-TextValue $ws.Range("E32") '  +2.02%  '

# Row 33
Set-TextValue $ws.Range("D33") '3.226'
Set-TextValue $ws.Range("E33") '  +1.11%  '

# Row 34
Set-TextValue $ws.Range("D34") '1.564'
Set-TextValue $ws.Range("E34") '  +1.16%  '

# Row 36
Set-TextValue $ws.Range("D36") '0.9049'
Set-TextValue $ws.Range("E36") '  -0.49%  '

# Row 37
Set-TextValue $ws.Range("D37") '2.581'
Set-TextValue $ws.Range("E37") '  +0.37%  '

# Row 38
$ws.Range("B38").Value = 'ImmutableX'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue $ws.Range("D38") '0.5535'
Set-TextValue $ws.Range("E38") '  +0.86%  '

# Row 39
$ws.Range("B39").Value = 'Maker'
$ws.Range("C39").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue $ws.Range("D39") '1.131.13'
Set-TextValue $ws.Range("E39") '  +0.75%  '

# Row 40
Set-TextValue $ws.Range("E40") '  +0.47%  '

# Row 41
Set-TextValue $ws.Range("D41") '0.9966'
Set-TextValue $ws.Range("E41") '  -0.43%  '

# Row 42
Set-TextValue $ws.Range("D42") '5.538'
Set-TextValue $ws.Range("E42") '  -0.64%  '

# Row 43
Set-TextValue $ws.Range("D43") '0.8015'
Set-TextValue $ws.Range("E43") '  -0.03%  '

# Row 44
Set-TextValue $ws.Range("D44") '97.95'
Set-TextValue $ws.Range("E44") '  -0.49%  '

# Row 45
Set-TextValue $ws.Range("D45") '1.775.61'
Set-TextValue $ws.Range("E45") '  +0.26%  '

# Row 46
Set-TextValue $ws.Range("E46") '  -8.82%  '

# Row 47
Set-TextValue $ws.Range("D47") '55.52'
Set-TextValue $ws.Range("E47") '  +1.11%  '

# Row 48
Set-TextValue $ws.Range("D48") '0.4369'
Set-TextValue $ws.Range("E48") '  -2.38%  '

# Row 49
Set-TextValue $ws.Range("D49") '7.700'

# Row 50
Set-TextValue $ws.Range("E50") '  -2.11%  '

# Row 51
Set-TextValue $ws.Range("D51") '1.005'
Set-TextValue $ws.Range("E51") '  +0.42%  '
